$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column widths (best-fit widths recomputed by Excel for the table
#    B:J after the data refresh). Values below are the ColumnWidth
#    inputs that reproduce the target best-fit pixel widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 16.08333333333121
$ws.Columns.Item(3).ColumnWidth = 9.91666666666788
$ws.Columns.Item(4).ColumnWidth = 30.91666666666697
$ws.Columns.Item(5).ColumnWidth = 11.91666666666606
$ws.Columns.Item(6).ColumnWidth = 8.583333333333485
$ws.Columns.Item(7).ColumnWidth = 12.583333333334394
$ws.Columns.Item(8).ColumnWidth = 17.083333333332575
$ws.Columns.Item(9).ColumnWidth = 15.91666666666697
$ws.Columns.Item(10).ColumnWidth = 13.416666666664696

# ---------------------------------------------------------------------
# 2) Re-anchor the logo picture: narrower B/C columns shift the image
#    left. Re-assert its original size since resizing the host columns
#    also resizes a "move and size with cells" picture.
# ---------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Width = 76.81889763779527
$shp.Height = 48.188976377952756
$shp.Left = 61.91281680610236
$shp.Top = 19.405511811023622

# ---------------------------------------------------------------------
# 3) New "parte 1" account-statement figures: period values swap
#    (2312/2311 -> 2311/2312) and the balance is updated for both rows.
# ---------------------------------------------------------------------
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2312"
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
